$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "rate_honor" / "100000" column (column D) entirely, shifting
# the columns to its right (catatan, nilai, tgl_mitra_diterima,
# tgl_ikut_survei) one position to the left.
$ws.Range("D1:D2").EntireColumn.Delete()

# Column C ("vol") ends up with the leftover duplicate style that used to be
# on D1 ("rate_honor"); reset it back to plain Text formatting with no fill
# (same as A1/B1) so it matches the original "vol" header styling instead of
# the orphaned duplicate style.
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Interior.Pattern = -4142
